# "switching A15 and M67Lab"
#
# The calendar table has two rows that need their "Activity/Lab" cell
# (column 3) and its matching colored marker cell (column 5) swapped:
#   Row 21 ("Fri 2/28")  currently: "Activity 15: Errors and Power"     / (blank marker)
#   Row 22 ("Mon 3/3")   currently: "Module 6 and 7 Lab: Arsenic"       / colored "Module 6 and 7 Lab" marker
#
# After the edit:
#   Row 21 ("Fri 2/28")  -> "Module 6 and 7 Lab: Arsenic"               / colored "Module 6 and 7 Lab" marker
#   Row 22 ("Mon 3/3")   -> "Activity 15: Errors and Power"             / (marker removed, blank)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Row 21, Col 3: "Activity 15: Errors and Power" -> "Module 6 and 7 Lab: Arsenic"
$row21Col3 = '<w:p w14:paraId="75582FED" w14:textId="368C0D5A" w:rsidR="000E022C" w:rsidRPr="009F5EDD" w:rsidRDefault="000E022C" w:rsidP="000E022C"><w:pPr><w:ind w:left="5"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Module 6 and 7 Lab: Arsenic</w:t></w:r></w:p>'
$t.Cell(21, 3).Range.InsertXML($pkgOpen + $row21Col3 + $pkgClose)

# --- Row 21, Col 5: blank marker -> colored "Module 6 and 7 Lab" marker (+ trailing blank paragraph)
$row21Col5 = '<w:p w14:paraId="78ED8895" w14:textId="7C3BDF96" w:rsidR="000E022C" w:rsidRPr="006B4DA2" w:rsidRDefault="000E022C" w:rsidP="000E022C"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="4472C4" w:themeColor="accent5"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="006B4DA2"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="4472C4" w:themeColor="accent5"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Module 6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="4472C4" w:themeColor="accent5"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="4472C4" w:themeColor="accent5"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">and 7 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="4472C4" w:themeColor="accent5"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Lab</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$t.Cell(21, 5).Range.InsertXML($pkgOpen + $row21Col5 + $pkgClose)

# --- Row 22, Col 3: "Module 6 and 7 Lab: Arsenic" -> "Activity 15: Errors and Power"
$row22Col3 = '<w:p w14:paraId="2B024943" w14:textId="4F208DB5" w:rsidR="000E022C" w:rsidRPr="006B4DA2" w:rsidRDefault="000D662E" w:rsidP="000E022C"><w:pPr><w:ind w:left="5"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Activity 15: Errors and Power</w:t></w:r></w:p>'
$t.Cell(22, 3).Range.InsertXML($pkgOpen + $row22Col3 + $pkgClose)

# --- Row 22, Col 5: colored "Module 6 and 7 Lab" marker removed, "Assignment 6" + blank paragraph remain
$row22Col5 = '<w:p w14:paraId="15649E71" w14:textId="77777777" w:rsidR="000E022C" w:rsidRDefault="000E022C" w:rsidP="000E022C"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="00B050"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="009F5EDD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="00B050"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Assignment </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="00B050"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>6</w:t></w:r></w:p><w:p w14:paraId="55CA64E2" w14:textId="7004265E" w:rsidR="000E022C" w:rsidRPr="00164261" w:rsidRDefault="000E022C" w:rsidP="000E022C"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$t.Cell(22, 5).Range.InsertXML($pkgOpen + $row22Col5 + $pkgClose)

Write-Host "Done"
